# Add three new columns (D, E, F) with header labels to Sheet1, matching
# the style already used by the existing header row (C1), and move the
# active selection to E6 (as captured in the saved workbook view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (C1, style index 1)
# onto the three new header cells so D1:F1 end up with the same style as
# A1:C1 instead of the worksheet default style.
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)  # xlPasteFormats

# New header values (these also introduce the three new shared strings:
# ORG_UN_IDENOLD, ORG_UN_IDENNEW, ORG_UN_STATUS).
$ws.Range("D1").Value = "ORG_UN_IDENOLD"
$ws.Range("E1").Value = "ORG_UN_IDENNEW"
$ws.Range("F1").Value = "ORG_UN_STATUS"

# Match the persisted selection/active-cell state from the edited workbook.
$ws.Range("E6").Select()
